$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.538.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.73%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.848.81"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +0.89%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.20%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'262.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +1.18%  "
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'  +0.10%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.5243"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +0.84%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.3230"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.06804"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +1.12%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'18.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +1.18%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.7827"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +2.22%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.07760"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +1.12%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'1.850.46"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.89%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'88.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -0.27%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'5.028"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +0.20%  "
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'  +0.02%  "
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = "'  -0.70%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'0.000007963"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.84%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  +0.06%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'26.586.52"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Value = "'4.643"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +2.56%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'9.459"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.50%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'5.993"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +1.44%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'143.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -1.24%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  -4.90%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  +2.21%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +0.71%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'112.12"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +0.77%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'4.184"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -0.17%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'0.08722"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -0.20%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'4.101"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -0.58%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +0.59%  "
$ws.Range("E32").ClearFormats()
$ws.Range("B33").Value = "'ARBITRUM"
$ws.Range("B33").ClearFormats()
$ws.Range("C33").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C33").ClearFormats()
$ws.Range("D33").Value = "'1.130"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +0.42%  "
$ws.Range("E33").ClearFormats()
$ws.Range("B34").Value = "'HuobiToken"
$ws.Range("B34").ClearFormats()
$ws.Range("C34").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C34").ClearFormats()
$ws.Range("D34").Value = "'2.881"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +1.28%  "
$ws.Range("E34").ClearFormats()
$ws.Range("B35").Value = "'ImmutableX"
$ws.Range("B35").ClearFormats()
$ws.Range("C35").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C35").ClearFormats()
$ws.Range("D35").Value = "'0.7187"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +5.45%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'3.104"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +0.60%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'2.276"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +2.85%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.01786"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +0.95%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.4854"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -0.92%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.9007"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'111.03"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -1.15%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'5.959"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -2.50%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  +0.10%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'7.668"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -0.43%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.4172"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -0.61%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.05877"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +0.09%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'8.984"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -0.55%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'35.14"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -0.49%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.1232"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -1.65%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.8931"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +3.54%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'59.97"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +1.25%  "
$ws.Range("E51").ClearFormats()
